# Display wrong verdict and violation level in results.
#
# The "Недочеты по политикам" (shortcomings by policy) table gains two new
# columns - "Ложные вердикты" (false verdicts) and "Ложные уровни нарушения"
# (false violation levels) - inserted just before the "Итого недочетов"
# (total shortcomings) column. The existing columns are also reordered so
# that "Ложные теги" (false tags) moves from the first data column to the
# last one before the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the two new columns (K, L) next to the existing ones:
#    same column width as the other data columns, and the same look as
#    the rest of the table (bordered, header cells centred/wrapped) -
#    achieved by copying the formatting from the neighbouring column.
# ---------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 14.17
$ws.Columns.Item(12).ColumnWidth = 14.17

$ws.Range("E3").Copy()
$ws.Range("K3:L3").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("K4:L16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Re-order / rewrite the header row (row 3).
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "Ложные политики"
$ws.Range("F3").Value = "Отсутствующие политики"
$ws.Range("G3").Value = "Ложные объекты"
$ws.Range("H3").Value = "Отсутствующие объекты"
$ws.Range("I3").Value = "Ложные теги"
$ws.Range("J3").Value = "Ложные вердикты"
$ws.Range("K3").Value = "Ложные уровни нарушения"
$ws.Range("L3").Value = "Итого недочетов"

# ---------------------------------------------------------------------
# 3. Rewrite the data rows (4-16) for columns E..L with the new values.
#    Column order is: E F G H I J K L
# ---------------------------------------------------------------------
$rows = @{
    4  = @(5, 0, 0, 0, 1, 0, 0, 6)
    5  = @(2, 0, 4, 0, 0, 0, 0, 6)
    6  = @(0, 0, 0, 0, 1, 0, 0, 1)
    7  = @(0, 0, 0, 1, 0, 0, 0, 1)
    8  = @(0, 2, 0, 0, 2, 0, 2, 6)
    9  = @(0, 0, 0, 0, 0, 0, 0, 0)
    10 = @(0, 0, 0, 0, 0, 0, 0, 0)
    11 = @(0, 0, 0, 0, 0, 0, 0, 0)
    12 = @(0, 0, 0, 0, 0, 0, 0, 0)
    13 = @(0, 4, 0, 0, 1, 0, 1, 6)
    14 = @(0, 1, 0, 0, 1, 0, 1, 3)
    15 = @(0, 1, 0, 1, 2, 1, 1, 6)
    16 = @(0, 0, 0, 0, 0, 1, 0, 1)
}

$cols = @("E", "F", "G", "H", "I", "J", "K", "L")

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# 4. Fix the two mismatched counters in the "object" table further down.
# ---------------------------------------------------------------------
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
